$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")
$ws.Range("B2").Value = 21
